$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add <w:outlineLvl w:val="0"/> to several paragraphs' pPr (OutlineLevel=1
#    in the Word object model maps to w:outlineLvl w:val="0").
# ---------------------------------------------------------------------------
$d.Paragraphs(1).Range.ParagraphFormat.OutlineLevel = 1   # "Yerba Buena, ..."
$d.Paragraphs(2).Range.ParagraphFormat.OutlineLevel = 1   # "ORDENANZA Nº 1890"
$d.Paragraphs(3).Range.ParagraphFormat.OutlineLevel = 1   # "VISTO:"
$d.Paragraphs(4).Range.ParagraphFormat.OutlineLevel = 1   # " La Ordenanza ..."
$d.Paragraphs(5).Range.ParagraphFormat.OutlineLevel = 1   # "CONSIDERANDO:"
$d.Paragraphs(10).Range.ParagraphFormat.OutlineLevel = 1  # "ARTÍCULO SEGUNDO: COMUNIQUESE..."

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the end of paragraph 7
#    ("EL CONCEJO DELIBERANTE...") to right at the start of paragraph 2's
#    content ("ORDENANZA Nº 1890").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$p2 = $d.Paragraphs(2)
$bmRange = $p2.Range
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 3) Remove the stray leading-space run at the start of paragraph 4
#    (" La Ordenanza Nº 1860 del 31/05/12; y").
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$p4Start = $p4.Range.Start
$d.Range($p4Start, $p4Start + 1).Delete()

# ---------------------------------------------------------------------------
# 4) Remove the stray leading-space run at the start of paragraph 6
#    (" Que el Departamento Ejecutivo Municipal propone cambios ...").
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs(6)
$p6Start = $p6.Range.Start
$d.Range($p6Start, $p6Start + 1).Delete()

# ---------------------------------------------------------------------------
# 5) Tighten "B) " -> "B)" (drop the trailing space after the closing
#    parenthesis) in paragraph 9.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$parenIdx = $full.IndexOf(") ")
if ($parenIdx -ge 0) {
    $d.Range($parenIdx + 1, $parenIdx + 2).Delete()
}

# ---------------------------------------------------------------------------
# 6) Paragraph 10: "ARTÍCULO SEGUNDO: COMUNIQUESE, REGISTRESE Y ARCHIVESE. "
#    becomes "ARTÍCULO TERCERO:" (underlined) + " " (plain) + the rest as-is.
# ---------------------------------------------------------------------------
$p10 = $d.Paragraphs(10)
$p10Search = $d.Range($p10.Range.Start, $p10.Range.End)
$p10Search.Find.Execute("ARTÍCULO SEGUNDO", $false, $false, $false, $false, $false, $true, 1, $false, "ARTÍCULO TERCERO", 2)

$p10b = $d.Paragraphs(10)
$p10Start = $p10b.Range.Start
$p10Text = $p10b.Range.Text
$colonIdx = $p10Text.IndexOf(":")
$colonAbs = $p10Start + $colonIdx

$underlineRange = $d.Range($p10Start, $colonAbs + 1)
$underlineRange.Font.Underline = 1

# ---------------------------------------------------------------------------
# 7) Delete the trailing paragraph that contains only "´".
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Delete()
